$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.686.42'
$ws.Range("E2").Value = '  -0.72%  '

$ws.Range("D3").Value = '3.946.37'
$ws.Range("E3").Value = '  -2.57%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.89%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.26'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +11.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.682'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.786'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.64%  '

$ws.Range("E10").Value = '  +5.81%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '55.87'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.95%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000326'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.59'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.94%  '

$ws.Range("D14").Value = '4.574.44'
$ws.Range("E14").Value = '  -2.76%  '

$ws.Range("D15").Value = '3.969.24'
$ws.Range("E15").Value = '  -2.18%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.03'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.81%  '

$ws.Range("E18").Value = '  -1.22%  '

$ws.Range("D19").Value = '72.632.25'
$ws.Range("E19").Value = '  -0.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.130'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '444.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.49%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.85'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '95.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.84%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.32%  '

$ws.Range("E28").Value = '  -1.23%  '

$ws.Range("E29").Value = '  -4.31%  '

$ws.Range("E30").Value = '  -3.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.89'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.89'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.64%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '49.68'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.37%  '

$ws.Range("E34").Value = '  -4.36%  '

$ws.Range("D35").Value = '0.0₃0993'
$ws.Range("E35").Value = '  +13.39%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '69.13'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.25%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '632.71'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.34%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.429'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.43'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.28%  '

$ws.Range("E40").Value = '  +0.09%  '

$ws.Range("E41").Value = '  -1.55%  '

$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("E43").Value = '  -3.56%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.62'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.96%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.18'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +41.76%  '

$ws.Range("E46").Value = '  -2.48%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.91%  '

$ws.Range("E48").Value = '  -0.88%  '

$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.83'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -16.78%  '

$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000285'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.31%  '

$ws.Range("D51").Value = '2.829.22'
$ws.Range("E51").Value = '  +1.45%  '
